# Add a new "Zubehör" (accessory) row to the Eigen_Glasdach sheet, between
# the existing "Tiefe" row and the final "Preis" row, and extend the price
# formula to include the new p_z accessory variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eigen_Glasdach")

# Push the old row 5 ("Preis" / "Gesamtpreis" / "Endpreis" / formula) down to
# row 6, opening up a fresh row 5 for the new accessory line.
$ws.Rows.Item(5).Insert()

# New row 5: Mehrfach | Zubehör  | p_z | Glashalter:170, Anschluß Wand:35
$ws.Range("A5").Value = "Mehrfach"
$ws.Range("B5").Value = "Zubehör "
$ws.Range("C5").Value = "p_z"

# Update the formula (now on row 6) to add the new accessory price term.
$ws.Range("E6").Value = "(p_l * p_t * p_glas) + p_z"

# Options list for the new accessory row.
$ws.Range("D5").Value = "Glashalter:170, Anschluß Wand:35"

# Match the author's final selection/view state.
$ws.Range("D5").Select() | Out-Null
